# Split the three "Programa" / "Programa" (EN) / "Bibliografia" paragraphs
# into multiple w:t runs separated by <w:br/> manual line breaks, matching
# the boundaries shown in the OOXML diff.

$d = $word.ActiveDocument

function Break-At($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find.Execute failed for: $findText"
    }
}

# --- Paragraph "Programa" (Portuguese) ---
Break-At "científicas.2. Projetos de Monografia" "científicas.^l2. Projetos de Monografia"
Break-At "de engenharia.3. Métodos de Pesquisa" "de engenharia.^l3. Métodos de Pesquisa"
Break-At "bibliográfica.4. Normas Aplicáveis" "bibliográfica.^l4. Normas Aplicáveis"
Break-At "relatórios técnicos.5. Pesquisa em Bases de Dados" "relatórios técnicos.^l5. Pesquisa em Bases de Dados"
Break-At "em engenharia.6. Organização de Referências Bibliográficas" "em engenharia.^l6. Organização de Referências Bibliográficas"

# --- Paragraph "Programa" (English, italic) ---
Break-At "analyses.2. Monograph and Technical Report Projects" "analyses.^l2. Monograph and Technical Report Projects"
Break-At "engineering projects.3. Research Methods" "engineering projects.^l3. Research Methods"
Break-At "bibliographic research.4. Standards Applicable" "bibliographic research.^l4. Standards Applicable"
Break-At "technical reports.5. Search in Bibliographic Databases" "technical reports.^l5. Search in Bibliographic Databases"
Break-At "engineering research projects.6. Organization of Bibliographic References" "engineering research projects.^l6. Organization of Bibliographic References"

# --- Paragraph "Bibliografia" ---
Break-At "São Paulo. 2005. GIL, A.C." "São Paulo. 2005. ^lGIL, A.C."
Break-At "São Paulo, 2010.MIGUEL, PAULO CAUCHICK." "São Paulo, 2010.^lMIGUEL, PAULO CAUCHICK."
Break-At "Elsevier Brasil, 2014.NASCIMENTO, L. P." "Elsevier Brasil, 2014.^lNASCIMENTO, L. P."
Break-At "Cengage Learning, 2012.SANTOS, C. R." "Cengage Learning, 2012.^lSANTOS, C. R."

Write-Output "done"
